$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.083.90"
$ws.Range("E2").Value = "  +5.45%  "

# Row 3
$ws.Range("D3").Value = "3.517.17"
$ws.Range("E3").Value = "  +3.19%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.20%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.95%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").Value = "3.529.12"
$ws.Range("E8").Value = "  +3.48%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.578"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.32%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.36%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.125"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.61%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.438"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.70%  "

# Row 13
$ws.Range("D13").Value = "4.143.09"
$ws.Range("E13").Value = "  +3.80%  "

# Row 14
$ws.Range("E14").Value = "  +0.40%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.93%  "

# Row 16
$ws.Range("E16").Value = "  +5.49%  "

# Row 17
$ws.Range("D17").Value = "67.001.99"
$ws.Range("E17").Value = "  +5.20%  "

# Row 18
$ws.Range("D18").Value = "3.534.33"
$ws.Range("E18").Value = "  +3.87%  "

# Row 19
$ws.Range("E19").Value = "  +3.37%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.80%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "395.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.18%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.25%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.46%  "

# Row 24
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.06%  "

# Row 25
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000127"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +12.04%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.531"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.10%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.69%  "

# Row 28
$ws.Range("E28").Value = "  +2.30%  "

# Row 29
$ws.Range("E29").Value = "  +0.08%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.35%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.15%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.29%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.18%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.78%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.09%  "

# Row 36
$ws.Range("E36").Value = "  +5.55%  "

# Row 37
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.66"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.72%  "

# Row 38
$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.911"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.52%  "

# Row 39
$ws.Range("E39").Value = "  +7.05%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0749"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.50%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.22%  "

# Row 42
$ws.Range("E42").Value = "  +5.72%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.16%  "

# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.826.36"
$ws.Range("E44").Value = "  +0.16%  "

# Row 45
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "26.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.30%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.35%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.04%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0315"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.10%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "353.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.74%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.66%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.66%  "
